$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Draft" label in F1 to "Drafting of manuscript"
$ws.Range("F1").Value = "Drafting of manuscript"

# Update selected/active cell to E1
$ws.Range("E1").Select()
